# Generate Report for Archive
#
# 1) Status text changes from "Ready for handoff" to "In Translation"
#    on every sheet that reports handoff/translation status:
#      - Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
#      - zh-cn sheet: column C (Status), rows 2-3
#      - de-de sheet: column C (Status), rows 2-3
#
# 2) The (now narrower) Status columns are resized:
#      - Overview columns E & F
#      - zh-cn column C
#      - de-de column C

$wb = $excel.ActiveWorkbook

$new = "In Translation"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $new
$wsOverview.Range("F2").Value = $new
$wsOverview.Range("E3").Value = $new
$wsOverview.Range("F3").Value = $new
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $new
$wsZhCn.Range("C3").Value = $new
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $new
$wsDeDe.Range("C3").Value = $new
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
